# Auto-generated Excel COM-interop script to apply odds/value updates
# for the 2025-02-14 FlashScore weekly games workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 2.55
$ws.Range("L8").Value = 3.5
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6
$ws.Range("O8").Value = 1.62
$ws.Range("P8").Value = 2.2
$ws.Range("S8").Value = 4.8
$ws.Range("W8").Value = 1.67
$ws.Range("X8").Value = 2.1
$ws.Range("AA8").Value = 6.5
$ws.Range("AC8").Value = 12
$ws.Range("AF8").Value = 51
$ws.Range("AO8").Value = 26

# Row 9
$ws.Range("M9").Value = 1.14
$ws.Range("N9").Value = 5.5

# Row 10
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = 1.08
$ws.Range("W10").Value = 1.78
$ws.Range("X10").Value = 2.03

# Row 15
$ws.Range("H15").Value = 3.4
$ws.Range("I15").Value = 3.2
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 2.1
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.33
$ws.Range("P15").Value = 3.25
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.75
$ws.Range("S15").Value = 2.85
$ws.Range("T15").Value = 1.41
$ws.Range("U15").Value = 3.75
$ws.Range("V15").Value = 1.25
$ws.Range("W15").Value = 1.44
$ws.Range("X15").Value = 2.63
$ws.Range("Y15").Value = 1.8
$ws.Range("Z15").Value = 1.91
$ws.Range("AA15").Value = 7.5
$ws.Range("AB15").Value = 10
$ws.Range("AC15").Value = 9.5
$ws.Range("AE15").Value = 19
$ws.Range("AF15").Value = 29
$ws.Range("AG15").Value = 9.5
$ws.Range("AJ15").Value = 51
$ws.Range("AK15").Value = 251
$ws.Range("AL15").Value = 9.5
$ws.Range("AM15").Value = 15
$ws.Range("AN15").Value = 12
$ws.Range("AR15").Value = 1.53
$ws.Range("AS15").Value = 2.47

# Row 18
$ws.Range("K18").Value = 1.95
$ws.Range("O18").Value = 1.44
$ws.Range("P18").Value = 2.63
$ws.Range("Q18").Value = 2.4
$ws.Range("R18").Value = 1.53
$ws.Range("U18").Value = 4.5
$ws.Range("V18").Value = 1.18
$ws.Range("W18").Value = 1.53
$ws.Range("X18").Value = 2.38
$ws.Range("Y18").Value = 2.1
$ws.Range("Z18").Value = 1.67
$ws.Range("AB18").Value = 9
$ws.Range("AC18").Value = 9.5
$ws.Range("AG18").Value = 6.5
$ws.Range("AI18").Value = 17
$ws.Range("AL18").Value = 9
$ws.Range("AQ18").Value = 41
$ws.Range("AR18").Value = 1.85
$ws.Range("AS18").Value = 2

# Row 23
$ws.Range("G23").Value = 4.75
$ws.Range("H23").Value = 3.8
$ws.Range("I23").Value = 1.67
$ws.Range("L23").Value = 2.3
$ws.Range("U23").Value = 3.75
$ws.Range("V23").Value = 1.25
$ws.Range("Y23").Value = 2.1
$ws.Range("Z23").Value = 1.67
$ws.Range("AC23").Value = 17
$ws.Range("AF23").Value = 51
$ws.Range("AG23").Value = 9
$ws.Range("AI23").Value = 21
$ws.Range("AM23").Value = 7
$ws.Range("AN23").Value = 9
$ws.Range("AO23").Value = 12
$ws.Range("AQ23").Value = 34

# Row 55
$ws.Range("M55").Value = 1.07
$ws.Range("N55").Value = 9
$ws.Range("O55").Value = 1.4
$ws.Range("P55").Value = 2.75
$ws.Range("Q55").Value = 2.25
$ws.Range("R55").Value = 1.62
$ws.Range("AK55").Value = 800
$ws.Range("AR55").Value = 1.78
$ws.Range("AS55").Value = 2.03

# Row 57
$ws.Range("G57").Value = 2.8
$ws.Range("H57").Value = 3.3
$ws.Range("I57").Value = 2.4
$ws.Range("J57").Value = 3.4
$ws.Range("K57").Value = 2.25
$ws.Range("L57").Value = 3
$ws.Range("M57").Value = 1.04
$ws.Range("N57").Value = 13
$ws.Range("O57").Value = 1.22
$ws.Range("P57").Value = 4
$ws.Range("Q57").Value = 1.75
$ws.Range("R57").Value = 2.05
$ws.Range("S57").Value = 2.3
$ws.Range("T57").Value = 1.62
$ws.Range("U57").Value = 2.75
$ws.Range("V57").Value = 1.4
$ws.Range("W57").Value = 1.33
$ws.Range("X57").Value = 3.25
$ws.Range("Y57").Value = 1.57
$ws.Range("Z57").Value = 2.25
$ws.Range("AB57").Value = 15
$ws.Range("AC57").Value = 11
$ws.Range("AD57").Value = 29
$ws.Range("AE57").Value = 21
$ws.Range("AF57").Value = 26
$ws.Range("AG57").Value = 13
$ws.Range("AH57").Value = 6.5
$ws.Range("AI57").Value = 12
$ws.Range("AK57").Value = 126
$ws.Range("AL57").Value = 10
$ws.Range("AM57").Value = 13
$ws.Range("AO57").Value = 23
$ws.Range("AQ57").Value = 23
$ws.Range("AR57").Value = 1.37
$ws.Range("AS57").Value = 3.15

# Row 74
$ws.Range("G74").Value = 2.3
$ws.Range("I74").Value = 3.1
$ws.Range("J74").Value = 3.4
$ws.Range("K74").Value = 1.8
$ws.Range("L74").Value = 4.33
$ws.Range("O74").Value = 1.67
$ws.Range("P74").Value = 2.1
$ws.Range("U74").Value = 7
$ws.Range("V74").Value = 1.1
$ws.Range("AA74").Value = 5.5
$ws.Range("AB74").Value = 9.5
$ws.Range("AD74").Value = 23
$ws.Range("AE74").Value = 26
$ws.Range("AL74").Value = 6.5
$ws.Range("AM74").Value = 13

# Row 75
$ws.Range("N75").Value = 10
$ws.Range("O75").Value = 1.36
$ws.Range("P75").Value = 3
$ws.Range("Q75").Value = 2.15
$ws.Range("R75").Value = 1.67
$ws.Range("U75").Value = 4
$ws.Range("V75").Value = 1.22
$ws.Range("W75").Value = 1.41
$ws.Range("X75").Value = 2.62
